$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.401.23"
$ws.Range("E2").Value = "  -0.58%  "

$ws.Range("D3").Value = "1.873.89"
$ws.Range("E3").Value = "  -0.24%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4714"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.84%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2884"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.91%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06482"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.73%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.93"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.79%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07789"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.69%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.873.89"
$ws.Range("E12").Value = "  -0.20%  "

$ws.Range("B13").Value = "Litecoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "96.15"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.67%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7211"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.68%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.144"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.93%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "284.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.98%  "

$ws.Range("D17").Value = "30.397.52"
$ws.Range("E17").Value = "  -1.03%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.03%  "

$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9999"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.03%  "

$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007493"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.50%  "

$ws.Range("D21").Value = "2.117.18"
$ws.Range("E21").Value = "  -0.28%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.264"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.44%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.253"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.63%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "163.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.18%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.055"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.42%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.68%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.879"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.64%  "

$ws.Range("E29").Value = "  -1.36%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09591"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.31%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.487"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.17%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.229"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.31%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.136"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.28%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04844"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.94%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.121"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.41%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6882"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.13%  "

$ws.Range("E37").Value = "  -0.16%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01896"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.35%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.816"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.84%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "75.56"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.96%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.204"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.09%  "

$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4221"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.78%  "

$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.928"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.95%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9994"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.14%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8283"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.40%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "100.84"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.38%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.746"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.48%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.979"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.05%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.84%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "898.92"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.00%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05727"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.26%  "
